# Update the "Förändrad" (Changed) date column (C) for rows 2 through 28
# from 45549 to 45550 (i.e. bump the date by one day), keeping the
# existing date number format/style intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45549) {
        $cell.Value2 = 45550
    }
}
